$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new mail-log entry as row 44 ---
$logs = $wb.Worksheets.Item("Logs")

$newRow = $logs.UsedRange.Rows.Count + 1

$logs.Cells.Item($newRow, 1).Value = "Uitnodiging voor netwerkevent"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Graag nodig ik u uit voor ons zakelijke netwerkevent volgende maand."
$logs.Cells.Item($newRow, 4).Value = "Samenwerking / Partnerverzoek"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 22:31:10"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# --- Logs sheet: extend the conditional-formatting ranges to include row 44 ---
$catFc = $logs.Range("D2:D43").FormatConditions.Item(1)
$catFc.ModifyAppliesToRange($logs.Range("D2:D$newRow"))

$answeredFc = $logs.Range("G2:G43").FormatConditions.Item(1)
$answeredFc.ModifyAppliesToRange($logs.Range("G2:G$newRow"))

# --- Dashboard sheet: bump the "Samenwerking / Partnerverzoek" count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 12
